$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (columns B..AA, i.e. columns 2..27) holds a 0-based header numbering
# that runs across the top of the grid. It is being renumbered to start one
# lower: B1 (which held 0) is cleared, and every other cell's number is
# decremented by one so the sequence 0..24 occupies C1..AA1 instead of
# 0..25 occupying B1..AA1.
for ($col = 2; $col -le 27; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $old = $cell.Value2
    if ($old -eq 0) {
        $cell.ClearContents()
    } else {
        $cell.Value = $old - 1
    }
}

# Column B (rows 2..26) holds the matching row index down the left side of
# the grid; renumber it the same way (decrement every value by one).
for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $old = $cell.Value2
    $cell.Value = $old - 1
}

# Move the active selection to reflect where the author was last working.
$ws.Range("M7").Select()
